$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.796.67'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '2.105.39'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''227.61'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('D7').Value = '''61.98'
$ws.Range('E7').Value = '  +2.48%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = '''16.03'
$ws.Range('E12').Value = '  +7.34%  '
$ws.Range('D13').Value = '2.417.07'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').Value = '''21.99'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '''5.49'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '2.065.41'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').Value = '38.967.33'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '''71.65'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '''6.04'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = '''227.23'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = '''2.31'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = '''9.66'
$ws.Range('E26').Value = '  +2.04%  '
$ws.Range('D27').Value = '''170.77'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('D30').Value = '''19.35'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('E31').Value = '  +9.55%  '
$ws.Range('D32').Value = '''0.121'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('B34').Value = 'THORChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D34').Value = '''7.15'
$ws.Range('E34').Value = '  +11.64%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''4.78'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '''3.51'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '''0.0230'
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('D41').Value = '''17.95'
$ws.Range('E41').Value = '  -1.97%  '
$ws.Range('D42').Value = '''102.07'
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('E44').Value = '  +7.41%  '
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('E47').Value = '  +1.80%  '
$ws.Range('E48').Value = '  +5.45%  '
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('D50').Value = '''2.96'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('D51').Value = '2.302.73'
$ws.Range('E51').Value = '  +0.68%  '
